$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Apparatus" (sheet2.xml): insert a new header row (row 3)
# with parameter names for the GFL inverter columns, shifting the
# existing rows 3-4 down to rows 4-5, and update row "6" (old row 5)
# values/formula.
# ---------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("Apparatus")

# Insert a new row above current row 3 (pushes rows 3-5 down to 4-6)
$wsApp.Rows.Item(3).Insert()

# New header row 3: labels for each parameter column (C..I)
$wsApp.Range("C3").Value = "Vdc (pu)"
$wsApp.Range("D3").Value = "Cdc (pu)"
$wsApp.Range("E3").Value = "wL (pu)"
$wsApp.Range("F3").Value = "R (pu)"
$wsApp.Range("G3").Value = "BW vdc (Hz)"
$wsApp.Range("H3").Value = "BW PLL (Hz)"
$wsApp.Range("I3").Value = "BW idq (Hz)"

# Row 6 (previously row 5) updates: F6 becomes a formula, I6 changes to 300
$wsApp.Range("F6").Formula = "=E6/10"
$wsApp.Range("I6").Value = 300

$wsApp.Range("G6").Select()

# ---------------------------------------------------------------
# Sheet "Advance" (sheet6.xml)
# ---------------------------------------------------------------
$wsAdv = $wb.Worksheets.Item("Advance")

$wsAdv.Range("B8").Value = 0

$wsAdv.Range("B9").Select()

# ---------------------------------------------------------------
# Sheet "NetworkLine" (sheet4.xml) - edited/selected last so it
# remains the active (tabSelected) sheet, matching the original file.
# ---------------------------------------------------------------
$wsNL = $wb.Worksheets.Item("NetworkLine")

$wsNL.Range("C11").Formula = "=D11/10"
$wsNL.Range("D11").Value = 0.1

# Clear out row 12 contents (A12:D12, F12:G12), keep E12 style but clear value
$wsNL.Range("A12:D12").ClearContents()
$wsNL.Range("E12").ClearContents()
$wsNL.Range("F12:G12").ClearContents()

$wsNL.Range("C12").Select()
